$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.863.16"
$ws.Range("E2").Value = "  -2.12%  "
$ws.Range("D3").Value = "1.816.34"
$ws.Range("E3").Value = "  -1.06%  "
$ws.Range("E4").Value = "  -0.47%  "
$ws.Range("E5").Value = "  -0.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "308.37"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4610"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3639"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.41%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07217"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8573"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.70"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.79%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07500"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.26%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.817.81"
$ws.Range("E13").Value = "  -3.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.317"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.494"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.52"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.77%  "
$ws.Range("E17").Value = "  -0.25%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008579"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.006"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.44%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.41"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.72%  "
$ws.Range("D21").Value = "26.749.85"
$ws.Range("E21").Value = "  -3.26%  "
$ws.Range("E22").Value = "  -3.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.49"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "2.095.57"
$ws.Range("E24").Value = "  -0.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.59"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.843"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.074"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.084"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.16"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.96%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08848"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.952"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.406"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.23%  "
$ws.Range("E34").Value = "  -4.39%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7128"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.83%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.073"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.86%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05230"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.414"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01914"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.921"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.125"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.83%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5133"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.73%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1619"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.149"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4782"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.62%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.007"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.44%  "
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.89"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.01%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.01"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.52%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06238"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.08%  "
$ws.Range("E50").Value = "  -3.84%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.99"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.72%  "
